$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 98.333336
$ws.Range("I11").Value = 98.333336
$ws.Range("K11").Value = 98.333336
$ws.Range("M11").Value = 41.666664

$ws.Range("H62").Value = 3999
$ws.Range("I62").Value = 1998.5
$ws.Range("K62").Value = 1998.5
$ws.Range("M62").Value = -1374.5

$ws.Range("H65").Value = 3999
$ws.Range("I65").Value = 1998.5
$ws.Range("K65").Value = 9992.5
$ws.Range("M65").Value = -6872.5

$ws.Range("H76").Value = 1774.25
$ws.Range("I76").Value = 1365.6666
$ws.Range("K76").Value = 1365.6666
$ws.Range("M76").Value = -1050.6666

$ws.Range("H79").Value = 1774.25
$ws.Range("I79").Value = 1365.6666
$ws.Range("K79").Value = 1365.6666
$ws.Range("M79").Value = -273.6666

$ws.Range("H100").Value = 2525.8
$ws.Range("I100").Value = 2907.25
$ws.Range("K100").Value = 2907.25
$ws.Range("M100").Value = -2366.25

$ws.Range("H131").Value = 2980
$ws.Range("I131").Value = 1200
$ws.Range("J131").Value = 3425
$ws.Range("K131").Value = 3600
$ws.Range("L131").Value = 10275
$ws.Range("M131").Value = 1440
$ws.Range("N131").Value = -20355

$ws.Range("H132").Value = 1098.8182
$ws.Range("I132").Value = 1098.8182
$ws.Range("K132").Value = 3296.4546
$ws.Range("M132").Value = -766.4546

$ws.Range("H138").Value = 4414.636
$ws.Range("I138").Value = 1274
$ws.Range("J138").Value = 6209.2856
$ws.Range("K138").Value = 3822
$ws.Range("L138").Value = 18627.8568
$ws.Range("M138").Value = 1318
$ws.Range("N138").Value = -28907.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1406
$ws.Range("I61").Value = 1406
$ws.Range("K61").Value = 1406
$ws.Range("M61").Value = -1194

$ws.Range("H74").Value = 21652.783
$ws.Range("I74").Value = 20818.863
$ws.Range("K74").Value = 20818.863
$ws.Range("M74").Value = -19944.863

$ws.Range("H77").Value = 21652.783
$ws.Range("I77").Value = 20818.863
$ws.Range("K77").Value = 104094.315
$ws.Range("M77").Value = -99726.315

$ws.Range("H136").Value = 1406
$ws.Range("I136").Value = 1406
$ws.Range("K136").Value = 4218
$ws.Range("M136").Value = -1668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3000
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 3000
$ws.Range("N99").Value = -5996
$ws.Range("M99").ClearContents()

$ws.Range("H107").Value = 971.1429000000001
$ws.Range("I107").Value = 959.6
$ws.Range("K107").Value = 959.6
$ws.Range("M107").Value = 960.4

$ws.Range("H134").Value = 3999.25
$ws.Range("I134").Value = 3999
$ws.Range("K134").Value = 11997
$ws.Range("M134").Value = -9462

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 6666745
$ws.Range("I6").Value = 20000000
$ws.Range("K6").Value = 20000000
$ws.Range("M6").Value = -19999887

$ws.Range("H7").Value = 78.125
$ws.Range("I7").Value = 70.666664
$ws.Range("J7").Value = 100.5
$ws.Range("K7").Value = 70.666664
$ws.Range("L7").Value = 100.5
$ws.Range("M7").Value = 42.333336
$ws.Range("N7").Value = -326.5

$ws.Range("H25").Value = 11002.75
$ws.Range("I25").Value = 1337
$ws.Range("J25").Value = 40000
$ws.Range("K25").Value = 1337
$ws.Range("L25").Value = 40000
$ws.Range("M25").Value = -1163
$ws.Range("N25").Value = -40348

$ws.Range("H132").Value = 1944.5
$ws.Range("I132").Value = 1342
$ws.Range("J132").Value = 2747.8333
$ws.Range("K132").Value = 4026
$ws.Range("L132").Value = 8243.499899999999
$ws.Range("M132").Value = -1496
$ws.Range("N132").Value = -13303.4999

$ws.Range("H134").Value = 5155.615
$ws.Range("I134").Value = 4274.909
$ws.Range("K134").Value = 12824.727
$ws.Range("M134").Value = -10289.727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5054.625
$ws.Range("I4").Value = 5763.143
$ws.Range("K4").Value = 17289.429
$ws.Range("M4").Value = -17177.429

$ws.Range("H6").Value = 457.14285
$ws.Range("I6").Value = 200
$ws.Range("J6").Value = 2000
$ws.Range("K6").Value = 600
$ws.Range("L6").Value = 6000
$ws.Range("M6").Value = -487
$ws.Range("N6").Value = -6226

$ws.Range("H10").Value = 928.75
$ws.Range("I10").Value = 928.75
$ws.Range("K10").Value = 2786.25
$ws.Range("M10").Value = -2647.25

$ws.Range("H13").Value = 1410.8572
$ws.Range("I13").Value = 1767
$ws.Range("K13").Value = 5301
$ws.Range("M13").Value = -5133

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 16667317
$ws.Range("I12").Value = 16667317
$ws.Range("K12").Value = 16667317
$ws.Range("M12").Value = -16667177

$ws.Range("H107").Value = 2531.8333
$ws.Range("I107").Value = 3524.25
$ws.Range("K107").Value = 3524.25
$ws.Range("M107").Value = -1604.25

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 1662.0526
$ws.Range("I132").Value = 1315.4375
$ws.Range("J132").Value = 3510.6667
$ws.Range("K132").Value = 3946.3125
$ws.Range("L132").Value = 10532.0001
$ws.Range("M132").Value = -1416.3125
$ws.Range("N132").Value = -15592.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1933
$ws.Range("I7").Value = 1899.5
$ws.Range("K7").Value = 1899.5
$ws.Range("M7").Value = -1787.5

$ws.Range("H11").Value = 1500
$ws.Range("I11").Value = 1500
$ws.Range("K11").Value = 1500
$ws.Range("M11").Value = -1360

$ws.Range("H14").Value = 2950
$ws.Range("I14").Value = 1100
$ws.Range("J14").Value = 4800
$ws.Range("K14").Value = 1100
$ws.Range("L14").Value = 4800
$ws.Range("M14").Value = -928
$ws.Range("N14").Value = -5144

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

$ws.Range("H126").Value = 1933
$ws.Range("I126").Value = 1899.5
$ws.Range("K126").Value = 5698.5
$ws.Range("M126").Value = -3228.5

$ws.Range("H132").Value = 2931.6667
$ws.Range("I132").Value = 2457.9333
$ws.Range("J132").Value = 3721.2222
$ws.Range("K132").Value = 7373.7999
$ws.Range("L132").Value = 11163.6666
$ws.Range("M132").Value = -4843.7999
$ws.Range("N132").Value = -16223.6666

$ws.Range("H136").Value = 4031.3076
$ws.Range("I136").Value = 4050.5454
$ws.Range("K136").Value = 12151.6362
$ws.Range("M136").Value = -9601.636200000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 3000000
$ws.Range("I7").Value = 3000000
$ws.Range("K7").Value = 3000000
$ws.Range("M7").Value = -2999887

$ws.Range("H126").Value = 1447.5714
$ws.Range("I126").Value = 1195.1428
$ws.Range("K126").Value = 3585.4284
$ws.Range("M126").Value = -1115.4284

$ws.Range("H132").Value = 1464.8462
$ws.Range("I132").Value = 738.2857
$ws.Range("J132").Value = 2312.5
$ws.Range("K132").Value = 2214.8571
$ws.Range("L132").Value = 6937.5
$ws.Range("M132").Value = 315.1428999999998
$ws.Range("N132").Value = -11997.5

$ws.Range("H136").Value = 1588.3334
$ws.Range("I136").Value = 1588.3334
$ws.Range("K136").Value = 4765.0002
$ws.Range("M136").Value = -2215.0002
